$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (Sale ID 1) with new values
$ws.Range("B2").Value = 7
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 225
$ws.Range("E2").Value = 45737.58293981481

# Update existing row 3 (Sale ID 2) with new values
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 56
$ws.Range("D3").Value = 3360
$ws.Range("E3").Value = 45737.58299768518

# Add new rows 4-9 (Sale IDs 3-8)
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 280
$ws.Range("E4").Value = 45737.58461805555

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 8
$ws.Range("C5").Value = 56
$ws.Range("D5").Value = 280
$ws.Range("E5").Value = 45737.60230324074

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 120
$ws.Range("E6").Value = 45737.61083333333

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 10
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 10
$ws.Range("E7").Value = 45737.62180555556

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 8
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 5
$ws.Range("E8").Value = 45737.97798611111

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 10
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 5
$ws.Range("E9").Value = 45738.01177083333

# Apply the same Timestamp number format used in column E to the new rows
$ws.Range("E4:E9").NumberFormat = $ws.Range("E2").NumberFormat
